$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New meeting-minutes row (10/14 meeting) with values for each column.
$ws.Cells.Item(9,1).Value = "October 14 2023"
$ws.Cells.Item(9,3).Value = "6:00PM"
$ws.Cells.Item(9,6).Value = "User sends a message -> sent to the server -> writes to room file"
$ws.Cells.Item(9,4).Value = "7:05PM"
$ws.Cells.Item(9,5).Value = "Redid use case diagram"
$ws.Cells.Item(9,2).Value = "David, Madison, Sean"

# Copy the formatting from the row above so the new row matches the table style.
$null = $ws.Range("A8:E8").Copy()
$null = $ws.Range("A9:E9").PasteSpecial(-4122)
$null = ($excel.CutCopyMode = $false)

# Clean up the old trailing blank row - only B10 keeps its formatting.
$null = $ws.Range("A10").Clear()
$null = $ws.Range("C10:K10").Clear()

# Widen column F to fit the new "What was Confirmed" text.
$ws.Columns.Item(6).ColumnWidth = 52.67

# Move the active selection as recorded when the sheet was last saved.
$null = $ws.Range("B23").Select()
